$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells to avoid Excel auto-converting
# them to actual numbers (which would lose trailing zeros / punctuation formatting).
$ws.Range("D2").Value = '69.613.91'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '3.745.60'
$ws.Range("E3").Value = '  +5.71%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '611.00'
$ws.Range("E5").Value = '  +3.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.46'
$ws.Range("E6").Value = '  -3.78%  '
$ws.Range("D7").Value = '3.742.74'
$ws.Range("E7").Value = '  +5.68%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("E10").Value = '  +5.05%  '
$ws.Range("E11").Value = '  -3.39%  '
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.66'
$ws.Range("E13").Value = '  +6.07%  '
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").Value = '4.371.87'
$ws.Range("E15").Value = '  +5.89%  '
$ws.Range("D16").Value = '3.746.01'
$ws.Range("E16").Value = '  +6.04%  '
$ws.Range("D17").Value = '69.648.36'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("E19").Value = '  +1.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '512.56'
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("E21").Value = '  -1.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.51'
$ws.Range("E22").Value = '  +4.84%  '
$ws.Range("E23").Value = '  -1.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.74'
$ws.Range("E24").Value = '  +1.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.50'
$ws.Range("E25").Value = '  +5.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.25'
$ws.Range("E26").Value = '  -0.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.02'
$ws.Range("E27").Value = '  +3.31%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000127'
$ws.Range("E29").Value = '  +16.39%  '
$ws.Range("E30").Value = '  -1.34%  '
$ws.Range("E31").Value = '  +4.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.79'
$ws.Range("E32").Value = '  -2.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.32'
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.115'
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.20'
$ws.Range("E36").Value = '  +1.99%  '
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("E38").Value = '  +2.11%  '
$ws.Range("E39").Value = '  +3.59%  '
$ws.Range("E40").Value = '  +3.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.15'
$ws.Range("E41").Value = '  +1.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '44.52'
$ws.Range("E42").Value = '  -4.05%  '
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '419.09'
$ws.Range("E44").Value = '  +5.22%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '3.068.33'
$ws.Range("E45").Value = '  +1.88%  '
$ws.Range("E46").Value = '  -2.72%  '
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.64'
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("E49").Value = '  +2.68%  '
$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '135.64'
$ws.Range("E51").Value = '  +0.68%  '
